$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1496.7167
$ws.Range("I15").Value = 1496.7167
$ws.Range("K15").Value = 4490.1501
$ws.Range("M15").Value = -4321.1501
$ws.Range("H51").Value = 2333.3333
$ws.Range("I51").Value = 1640
$ws.Range("J51").Value = 3200
$ws.Range("K51").Value = 1640
$ws.Range("L51").Value = 3200
$ws.Range("M51").Value = -1156
$ws.Range("N51").Value = -4168
$ws.Range("H55").Value = 137.5
$ws.Range("I55").Value = 137.5
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 137.5
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 76.5
$ws.Range("N55").ClearContents()
$ws.Range("H62").Value = 2386.8235
$ws.Range("I62").Value = 2054.889
$ws.Range("J62").Value = 2760.25
$ws.Range("K62").Value = 2054.889
$ws.Range("L62").Value = 2760.25
$ws.Range("M62").Value = -1430.889
$ws.Range("N62").Value = -4008.25
$ws.Range("H65").Value = 2386.8235
$ws.Range("I65").Value = 2054.889
$ws.Range("J65").Value = 2760.25
$ws.Range("K65").Value = 10274.445
$ws.Range("L65").Value = 13801.25
$ws.Range("M65").Value = -7154.445
$ws.Range("N65").Value = -20041.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1367.1111
$ws.Range("I2").Value = 1038.2
$ws.Range("J2").Value = 1778.25
$ws.Range("K2").Value = 1038.2
$ws.Range("L2").Value = 1778.25
$ws.Range("M2").Value = -925.2
$ws.Range("N2").Value = -2004.25
$ws.Range("H45").Value = 1012
$ws.Range("I45").Value = 847
$ws.Range("J45").Value = 1507
$ws.Range("K45").Value = 847
$ws.Range("L45").Value = 1507
$ws.Range("M45").Value = -470
$ws.Range("N45").Value = -2261
$ws.Range("H116").Value = 1367.1111
$ws.Range("I116").Value = 1038.2
$ws.Range("J116").Value = 1778.25
$ws.Range("K116").Value = 1038.2
$ws.Range("L116").Value = 1778.25
$ws.Range("M116").Value = 1255.8
$ws.Range("N116").Value = -6366.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1367.1111
$ws.Range("I3").Value = 1038.2
$ws.Range("J3").Value = 1778.25
$ws.Range("K3").Value = 1038.2
$ws.Range("L3").Value = 1778.25
$ws.Range("M3").Value = -924.2
$ws.Range("N3").Value = -2006.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 900
$ws.Range("J6").Value = 900
$ws.Range("L6").Value = 900
$ws.Range("N6").Value = -1126
$ws.Range("H7").Value = 86.411766
$ws.Range("I7").Value = 47.615383
$ws.Range("J7").Value = 212.5
$ws.Range("K7").Value = 47.615383
$ws.Range("L7").Value = 212.5
$ws.Range("M7").Value = 65.38461699999999
$ws.Range("N7").Value = -438.5
$ws.Range("H17").Value = 34333.332
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2826
$ws.Range("H22").Value = 281.66666
$ws.Range("I22").Value = 219.05556
$ws.Range("J22").Value = 657.3333
$ws.Range("K22").Value = 219.05556
$ws.Range("L22").Value = 657.3333
$ws.Range("M22").Value = 130.94444
$ws.Range("N22").Value = -1357.3333
$ws.Range("H41").Value = 3000
$ws.Range("I41").Value = 3000
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3000
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2572
$ws.Range("N41").ClearContents()
$ws.Range("H50").Value = 29000
$ws.Range("J50").Value = 29000
$ws.Range("L50").Value = 29000
$ws.Range("N50").Value = -30250
$ws.Range("H51").Value = 20491.666
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 25737.5
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 25737.5
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -27209.5
$ws.Range("H59").Value = 40110
$ws.Range("I59").Value = 40110
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 40110
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -38965
$ws.Range("N59").ClearContents()
$ws.Range("H60").Value = 35149.7
$ws.Range("J60").Value = 37887.125
$ws.Range("L60").Value = 37887.125
$ws.Range("N60").Value = -38909.125
$ws.Range("H61").Value = 20491.666
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 25737.5
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 25737.5
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -26433.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1439.7142
$ws.Range("I5").Value = 1632.5294
$ws.Range("J5").Value = 1141.7273
$ws.Range("K5").Value = 4897.5882
$ws.Range("L5").Value = 3425.1819
$ws.Range("M5").Value = -4785.5882
$ws.Range("N5").Value = -3649.1819
$ws.Range("H12").Value = 715737.7
$ws.Range("I12").Value = 90.166664
$ws.Range("J12").Value = 920208.4
$ws.Range("K12").Value = 270.499992
$ws.Range("L12").Value = 2760625.2
$ws.Range("M12").Value = -97.49999200000002
$ws.Range("N12").Value = -2760971.2
$ws.Range("H131").Value = 22223882
$ws.Range("I131").Value = 304
$ws.Range("J131").Value = 25001830
$ws.Range("K131").Value = 912
$ws.Range("L131").Value = 75005490
$ws.Range("M131").Value = 4128
$ws.Range("N131").Value = -75015570
$ws.Range("H135").Value = 1439.7142
$ws.Range("I135").Value = 1632.5294
$ws.Range("J135").Value = 1141.7273
$ws.Range("K135").Value = 14692.7646
$ws.Range("L135").Value = 10275.5457
$ws.Range("M135").Value = -12157.7646
$ws.Range("N135").Value = -15345.5457
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 745.6316
$ws.Range("J22").Value = 1107.5555
$ws.Range("L22").Value = 1107.5555
$ws.Range("N22").Value = -1697.5555
$ws.Range("H27").Value = 745.6316
$ws.Range("J27").Value = 1107.5555
$ws.Range("L27").Value = 1107.5555
$ws.Range("N27").Value = -1321.5555
$ws.Range("H40").Value = 3606.5833
$ws.Range("I40").Value = 3556.8
$ws.Range("J40").Value = 3642.1428
$ws.Range("K40").Value = 3556.8
$ws.Range("L40").Value = 3642.1428
$ws.Range("M40").Value = -3420.8
$ws.Range("N40").Value = -3914.1428
$ws.Range("H122").Value = 7816465.5
$ws.Range("I122").Value = 10003260
$ws.Range("J122").Value = 6485.4287
$ws.Range("K122").Value = 30009780
$ws.Range("L122").Value = 19456.2861
$ws.Range("M122").Value = -30007330
$ws.Range("N122").Value = -24356.2861
$ws.Range("H136").Value = 3503.2354
$ws.Range("I136").Value = 2734.6155
$ws.Range("J136").Value = 6001.25
$ws.Range("K136").Value = 8203.8465
$ws.Range("L136").Value = 18003.75
$ws.Range("M136").Value = -5653.8465
$ws.Range("N136").Value = -23103.75
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 50333.332
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 50333.332
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 50333.332
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -50963.332
$ws.Range("H79").Value = 50333.332
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 50333.332
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 50333.332
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -52517.332
$ws.Range("H81").Value = 125709.5
$ws.Range("I81").Value = 125709.5
$ws.Range("K81").Value = 251419
$ws.Range("M81").Value = -250358
$ws.Range("H84").Value = 125709.5
$ws.Range("I84").Value = 125709.5
$ws.Range("K84").Value = 1257095
$ws.Range("M84").Value = -1251791
$ws.Range("H87").Value = 143993.75
$ws.Range("H88").Value = 28089
$ws.Range("J88").Value = 28089
$ws.Range("L88").Value = 28089
$ws.Range("N88").Value = -28901
$ws.Range("H90").Value = 143993.75
$ws.Range("H91").Value = 28089
$ws.Range("J91").Value = 28089
$ws.Range("L91").Value = 28089
$ws.Range("N91").Value = -30897
$ws.Range("H113").Value = 922
$ws.Range("I113").Value = 893.55554
$ws.Range("J113").Value = 1050
$ws.Range("K113").Value = 2680.66662
$ws.Range("L113").Value = 3150
$ws.Range("M113").Value = -510.66662
$ws.Range("N113").Value = -7490
$ws.Range("H122").Value = 7577969.5
$ws.Range("I122").Value = 9261084
$ws.Range("K122").Value = 27783252
$ws.Range("M122").Value = -27780802
$ws.Range("H132").Value = 1948.2593
$ws.Range("I132").Value = 1121.9445
$ws.Range("J132").Value = 3600.889
$ws.Range("K132").Value = 3365.8335
$ws.Range("L132").Value = 10802.667
$ws.Range("M132").Value = -835.8335000000002
$ws.Range("N132").Value = -15862.667
$ws.Range("H136").Value = 2400.1538
$ws.Range("I136").Value = 2445.7273
$ws.Range("K136").Value = 7337.1819
$ws.Range("M136").Value = -4787.1819
